$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New week column H: header date 2022-10-17 (serial 44851), same date format as D1:G1
$ws.Range("H1").Value = 44851
$ws.Range("H1").NumberFormat = $ws.Range("G1").NumberFormat

# Fill in the new participation counts recorded for this date, plus a couple of
# values that were backfilled into the existing 3rd-Oct week (column F) and
# the very first week (column D/E/F) for row 2.
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 5
$ws.Range("H8").Value = 10
$ws.Range("H12").Value = 10
$ws.Range("H15").Value = 7
$ws.Range("H26").Value = 9
$ws.Range("H27").Value = 8
$ws.Range("H33").Value = 9
$ws.Range("H35").Value = 6
$ws.Range("H36").Value = 1

# Student rename: "Vladislav Taltos" -> "Vladislav Gorbachev"
$ws.Range("B33").Value = "Vladislav Gorbachev"

# Match the author's final selection
$ws.Range("H36").Select()
